# VCE scenario policy file updates
# Updates the "Qualifies for RPS (Boolean)" flags on the RQSD-BRQSD and
# RQSD-RQSD sheets, and leaves the workbook's active sheet/selection the
# way the author left it when saving (RQSD-BRQSD tab, cell B22 selected
# on both RQSD-BRQSD and RQSD-RQSD).

$wb = $excel.ActiveWorkbook

# --- RQSD-BRQSD (sheet2): "municipal solid waste" now qualifies ---
$wsBRQSD = $wb.Worksheets.Item("RQSD-BRQSD")
$wsBRQSD.Range("B17").Value = 1

# --- RQSD-RQSD (sheet3): nuclear, hydro, biomass, and municipal solid
#     waste now qualify ---
$wsRQSD = $wb.Worksheets.Item("RQSD-RQSD")
$wsRQSD.Range("B4").Value = 1
$wsRQSD.Range("B5").Value = 1
$wsRQSD.Range("B9").Value = 1
$wsRQSD.Range("B17").Value = 1

# --- column width touch-ups (as left by the author's Excel session) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Columns.Item(2).ColumnWidth = 83.5834125

$wsBRQSD.Columns.Item(1).ColumnWidth = 23.5834125
$wsBRQSD.Columns.Item(2).ColumnWidth = 29.250025

$wsRQSD.Columns.Item(1).ColumnWidth = 23.5834125
$wsRQSD.Columns.Item(2).ColumnWidth = 29.5834125

# --- selection / active sheet state ---
# RQSD-RQSD: leave cell B22 selected (not the active tab)
$wsRQSD.Activate()
$wsRQSD.Range("B22").Select()

# RQSD-BRQSD becomes the active/selected tab, with B22 selected
$wsBRQSD.Activate()
$wsBRQSD.Range("B22").Select()
